$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Swap the contents of row 18 and row 19 (columns F..V), the
#    "Indice" values in column A (and A..D which are identical anyway)
#    stay put - only the match data actually changes place.
# ------------------------------------------------------------------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($col in $cols) {
    $addr18 = "$col" + "18"
    $addr19 = "$col" + "19"
    $v18 = $ws.Range($addr18).Value2
    $v19 = $ws.Range($addr19).Value2
    $ws.Range($addr18).Value2 = $v19
    $ws.Range($addr19).Value2 = $v18
}

# ------------------------------------------------------------------
# 2) Append four new match rows (38-41), copying the formatting from
#    the last existing data row (37) so the bold/border style on
#    column A and the date style on column E are preserved.
# ------------------------------------------------------------------
$ws.Range("A37:V37").Copy() | Out-Null
$ws.Range("A38:V41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 38 (Indice 37): Sparta Selemet 0 - 1 Zimbru Chisinau
$ws.Range("A38").Value2 = 37
$ws.Range("B38").Value2 = "moldova"
$ws.Range("C38").Value2 = "super-liga"
$ws.Range("D38").Value2 = "2023-2024"
$ws.Range("E38").Value2 = 45227.54166666666
$ws.Range("F38").Value2 = "Sparta Selemet"
$ws.Range("G38").Value2 = 0
$ws.Range("H38").Value2 = "Zimbru Chisinau"
$ws.Range("I38").Value2 = 1
$ws.Range("J38").Value2 = 6.72
$ws.Range("K38").Value2 = "27/10/2023 01:12"
$ws.Range("L38").Value2 = 15.06
$ws.Range("M38").Value2 = "28/10/2023 12:59"
$ws.Range("N38").Value2 = 4.41
$ws.Range("O38").Value2 = "27/10/2023 01:12"
$ws.Range("P38").Value2 = 6.49
$ws.Range("Q38").Value2 = "28/10/2023 12:59"
$ws.Range("R38").Value2 = 1.29
$ws.Range("S38").Value2 = "27/10/2023 01:12"
$ws.Range("T38").Value2 = 1.13
$ws.Range("U38").Value2 = "28/10/2023 12:59"
$ws.Range("V38").Value2 = "https://www.betexplorer.com/football/moldova/super-liga/sparta-selemet-zimbru-chisinau/2s8TVQQL/"

# Row 39 (Indice 38): Milsami 1 - 1 Petrocub
$ws.Range("A39").Value2 = 38
$ws.Range("B39").Value2 = "moldova"
$ws.Range("C39").Value2 = "super-liga"
$ws.Range("D39").Value2 = "2023-2024"
$ws.Range("E39").Value2 = 45227.66666666666
$ws.Range("F39").Value2 = "Milsami"
$ws.Range("G39").Value2 = 1
$ws.Range("H39").Value2 = "Petrocub"
$ws.Range("I39").Value2 = 1
$ws.Range("J39").Value2 = 2.71
$ws.Range("K39").Value2 = "27/10/2023 04:12"
$ws.Range("L39").Value2 = 3.3
$ws.Range("M39").Value2 = "28/10/2023 15:41"
$ws.Range("N39").Value2 = 2.75
$ws.Range("O39").Value2 = "27/10/2023 04:12"
$ws.Range("P39").Value2 = 2.52
$ws.Range("Q39").Value2 = "28/10/2023 15:41"
$ws.Range("R39").Value2 = 2.42
$ws.Range("S39").Value2 = "27/10/2023 04:12"
$ws.Range("T39").Value2 = 2.47
$ws.Range("U39").Value2 = "28/10/2023 15:41"
$ws.Range("V39").Value2 = "https://www.betexplorer.com/football/moldova/super-liga/milsami-petrocub-hincesti/6FEKXnd9/"

# Row 40 (Indice 39): Dacia Buiucani 0 - 4 Sheriff Tiraspol
$ws.Range("A40").Value2 = 39
$ws.Range("B40").Value2 = "moldova"
$ws.Range("C40").Value2 = "super-liga"
$ws.Range("D40").Value2 = "2023-2024"
$ws.Range("E40").Value2 = 45228.54166666666
$ws.Range("F40").Value2 = "Dacia Buiucani"
$ws.Range("G40").Value2 = 0
$ws.Range("H40").Value2 = "Sheriff Tiraspol"
$ws.Range("I40").Value2 = 4
$ws.Range("J40").Value2 = 11
$ws.Range("K40").Value2 = "28/10/2023 02:13"
$ws.Range("L40").Value2 = 16.58
$ws.Range("M40").Value2 = "29/10/2023 12:44"
$ws.Range("N40").Value2 = 7.22
$ws.Range("O40").Value2 = "28/10/2023 02:13"
$ws.Range("P40").Value2 = 7.93
$ws.Range("Q40").Value2 = "29/10/2023 12:48"
$ws.Range("R40").Value2 = 1.1
$ws.Range("S40").Value2 = "28/10/2023 02:13"
$ws.Range("T40").Value2 = 1.1
$ws.Range("U40").Value2 = "29/10/2023 12:44"
$ws.Range("V40").Value2 = "https://www.betexplorer.com/football/moldova/super-liga/dacia-buiucani-sheriff-tiraspol/KIIOW6BF/"

# Row 41 (Indice 40): Floresti 2 - 3 Balti
$ws.Range("A41").Value2 = 40
$ws.Range("B41").Value2 = "moldova"
$ws.Range("C41").Value2 = "super-liga"
$ws.Range("D41").Value2 = "2023-2024"
$ws.Range("E41").Value2 = 45228.54166666666
$ws.Range("F41").Value2 = "Floresti"
$ws.Range("G41").Value2 = 2
$ws.Range("H41").Value2 = "Balti"
$ws.Range("I41").Value2 = 3
$ws.Range("J41").Value2 = 2.76
$ws.Range("K41").Value2 = "28/10/2023 02:13"
$ws.Range("L41").Value2 = 3.09
$ws.Range("M41").Value2 = "29/10/2023 12:58"
$ws.Range("N41").Value2 = 3.02
$ws.Range("O41").Value2 = "28/10/2023 02:13"
$ws.Range("P41").Value2 = 3.81
$ws.Range("Q41").Value2 = "29/10/2023 12:58"
$ws.Range("R41").Value2 = 2.2
$ws.Range("S41").Value2 = "28/10/2023 02:13"
$ws.Range("T41").Value2 = 1.93
$ws.Range("U41").Value2 = "29/10/2023 12:58"
$ws.Range("V41").Value2 = "https://www.betexplorer.com/football/moldova/super-liga/floresti-csf-balti/tY6XUptS/"

Write-Host "Edit complete"
